$wb = $excel.ActiveWorkbook

# Update the "Status" text in all sheets: "Ready for handoff" -> "In Translation"
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ("Ready for handoff" -eq $v) {
            $cell.Value2 = "In Translation"
        }
    }
}

# Adjust the "Status" column widths (auto-fit narrower to match the new text) on all three sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
